$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CMIP6 specialization id mappings for the "Top Of Atmos Insolation" /
# "Volcanoes" rows (these had a blank column C before).
# Values are entered in the same order the author appears to have typed
# them (this controls shared-string insertion order).

$ws.Range("C6").Value = "cmip6.atmos.solar/insolation_ozone.solar_ozone_impact"

$ws.Range("C7").Value = "cmip6.atmos.solar.orbital_parameters.computation_method"
$ws.Range("C8").Value = "cmip6.atmos.solar.orbital_parameters.computation_method"

$ws.Range("C11").Value = "cmip6.atmos.solar.orbital_parameters.fixed_reference_date"

$ws.Range("C13").Value = "cmip6.atmos.solar.orbital_parameters.type"
$ws.Range("C13").Borders.LineStyle = 1
$ws.Range("C13").Borders.Color = 10855845

$ws.Range("C15").Value = "cmip6.atmos.solar.solar_constant.fixed_value"
$ws.Range("C18").Value = "cmip6.atmos.solar.solar_constant.fixed_value"

$ws.Range("C12").Value = "cmip6.atmos.solar.orbital_parameters.transient_method"

$ws.Range("C16").Value = "cmip6.atmos.solar.solar_constant.transient_characteristics"

$ws.Range("C17").Value = "cmip6.atmos.solar.solar_constant.type"

$ws.Range("C19").Value = "cmip6.atmos.volcanos.volcanoes_treatment.volcanoes_implementation"

$ws.Range("C20").Select()
